# Add Week 15 simulations: new RB player "J.Hardy" row in the RB sheet.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("RB")

# New row 6 for the RB sheet: J.Hardy with all stats at 0
$ws.Range("A6").Value = "J.Hardy"
$ws.Range("B6:J6").Value = 0

# Update selection on the RB sheet and make it the active sheet/tab
$ws.Activate() | Out-Null
$ws.Range("L5").Select() | Out-Null
